# Swap the "category" and "group" columns in the SectorGroup worksheet.
#
# Column layout (A:G) = code, name, status, <D>, <E>, <F>, <G>
# Before this edit:  D = codeforiati:category-name, E = codeforiati:group-name,
#                     F = codeforiati:group-code,    G = codeforiati:category-code
# After this edit:   D = codeforiati:group-name,     E = codeforiati:category-name,
#                     F = codeforiati:category-code,  G = codeforiati:group-code
#
# i.e. for every row (including the header), the values in columns D and E
# are swapped with each other, and the values in columns F and G are
# swapped with each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Make sure columns D:G keep being stored as text (some of the values, e.g.
# "110", "111", look numeric and would otherwise be auto-converted to
# numbers by Excel).
$ws.Range("D1:G$lastRow").NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2

    $dCell.Value = $eVal
    $eCell.Value = $dVal
    $fCell.Value = $gVal
    $gCell.Value = $fVal
}
